$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: build a minimal pkg:package InsertXML payload for a single <w:p>
# body fragment (already-built inner XML for the paragraph's children).
# ---------------------------------------------------------------------------
function New-ParaPackage([string]$paraInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $paraInnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$rPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$pPr = '<w:pPr><w:pStyle w:val="PlainText"/>' + $rPr + '</w:pPr>'

function Find-ParagraphByText([string]$target) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($txt -eq $target) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "/ip firewall address-list remove [find where list="CountryIPBlocks"];"
#    -> split into 5 runs, with proofErr spellStart/spellEnd wrapping
#       "ip" and "CountryIPBlocks".
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphByText('/ip firewall address-list remove [find where list="CountryIPBlocks"];')
$inner1 = $pPr + `
    '<w:r>' + $rPr + '<w:t>/</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>ip</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> firewall address-list remove [find where list="</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>CountryIPBlocks</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>"];</w:t></w:r>'
$p1.Range.InsertXML((New-ParaPackage $inner1))

# ---------------------------------------------------------------------------
# 2) "/import file-name=IP-Firewall-Address-List.rsc"
#    -> split into 2 runs, proofErr spellStart/spellEnd wraps "List.rsc".
# ---------------------------------------------------------------------------
$p5 = Find-ParagraphByText('/import file-name=IP-Firewall-Address-List.rsc')
$inner5 = $pPr + `
    '<w:r>' + $rPr + '<w:t>/import file-name=IP-Firewall-Address-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>List.rsc</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$p5.Range.InsertXML((New-ParaPackage $inner5))

# ---------------------------------------------------------------------------
# 3) "#/file remove IP-Firewall-Address-List.rsc"
#    -> split into 2 runs, proofErr spellStart/spellEnd wraps "List.rsc".
# ---------------------------------------------------------------------------
$p6 = Find-ParagraphByText('#/file remove IP-Firewall-Address-List.rsc')
$inner6 = $pPr + `
    '<w:r>' + $rPr + '<w:t>#/file remove IP-Firewall-Address-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>List.rsc</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$p6.Range.InsertXML((New-ParaPackage $inner6))

# ---------------------------------------------------------------------------
# 4) Remove the "#AFGHANISTAN", "#PAKISTAN", "#ROMANIA", "#SERBIA" paragraphs
#    entirely (including their paragraph marks).
# ---------------------------------------------------------------------------
$toDelete = @("#AFGHANISTAN", "#PAKISTAN", "#ROMANIA", "#SERBIA")
foreach ($target in $toDelete) {
    $victim = Find-ParagraphByText($target)
    if ($victim -ne $null) {
        $victim.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 5) After "#TURKEY", insert two new paragraphs: "#TURKMENISTAN" (2 runs,
#    "#" + "TURKMENISTAN", no proofErr) and "#TURKS" (single run).
# ---------------------------------------------------------------------------
$pTurkey = Find-ParagraphByText('#TURKEY')
$pTurkey.Range.InsertParagraphAfter()
$newPara1 = Find-ParagraphByText('')
$innerTm = $pPr + `
    '<w:r>' + $rPr + '<w:t>#</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>TURKMENISTAN</w:t></w:r>'
$newPara1.Range.InsertXML((New-ParaPackage $innerTm))

$pTurkmenistan = Find-ParagraphByText('#TURKMENISTAN')
$pTurkmenistan.Range.InsertParagraphAfter()
$newPara2 = Find-ParagraphByText('')
$innerTk = $pPr + '<w:r>' + $rPr + '<w:t>#TURKS</w:t></w:r>'
$newPara2.Range.InsertXML((New-ParaPackage $innerTk))

Write-Output "Done"
